# Auto-generated Excel COM-interop script
# Fills in short-sale data (columns D-H) for rows 210-211 that were
# placeholder/blank, and appends new rows 212-217 (with 216-217 leaving
# columns D-H blank, matching the pre-existing placeholder pattern),
# across all 6 worksheets.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1 ----
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(210, 4).Value = 378845
$ws.Cells.Item(210, 5).Value = 204757766
$ws.Cells.Item(210, 6).Value = 101719882500
$ws.Cells.Item(210, 7).Value = 54977460171000
$ws.Cells.Item(210, 8).Value = 0.1899999976158142

$ws.Cells.Item(211, 4).Value = 403577
$ws.Cells.Item(211, 5).Value = 204757766
$ws.Cells.Item(211, 6).Value = 108562213000
$ws.Cells.Item(211, 7).Value = 55079839054000
$ws.Cells.Item(211, 8).Value = 0.2000000029802322

$ws.Cells.Item(212, 1).Value = 45968
$ws.Cells.Item(212, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 2).Value = 264000
$ws.Cells.Item(212, 3).Value = 666772
$ws.Cells.Item(212, 4).Value = 403079
$ws.Cells.Item(212, 5).Value = 204757766
$ws.Cells.Item(212, 6).Value = 106412856000
$ws.Cells.Item(212, 7).Value = 54056050224000
$ws.Cells.Item(212, 8).Value = 0.2000000029802322

$ws.Cells.Item(213, 1).Value = 45971
$ws.Cells.Item(213, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213, 2).Value = 270500
$ws.Cells.Item(213, 3).Value = 692464
$ws.Cells.Item(213, 4).Value = 405235
$ws.Cells.Item(213, 5).Value = 204757766
$ws.Cells.Item(213, 6).Value = 109616067500
$ws.Cells.Item(213, 7).Value = 55386975703000
$ws.Cells.Item(213, 8).Value = 0.2000000029802322

$ws.Cells.Item(214, 1).Value = 45972
$ws.Cells.Item(214, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214, 2).Value = 269000
$ws.Cells.Item(214, 3).Value = 494203
$ws.Cells.Item(214, 4).Value = 404845
$ws.Cells.Item(214, 5).Value = 204757766
$ws.Cells.Item(214, 6).Value = 108903305000
$ws.Cells.Item(214, 7).Value = 55079839054000
$ws.Cells.Item(214, 8).Value = 0.2000000029802322

$ws.Cells.Item(215, 1).Value = 45973
$ws.Cells.Item(215, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 2).Value = 275500
$ws.Cells.Item(215, 3).Value = 807526
$ws.Cells.Item(215, 4).Value = 453088
$ws.Cells.Item(215, 5).Value = 204757766
$ws.Cells.Item(215, 6).Value = 124825744000
$ws.Cells.Item(215, 7).Value = 56410764533000
$ws.Cells.Item(215, 8).Value = 0.2199999988079071

$ws.Cells.Item(216, 1).Value = 45974
$ws.Cells.Item(216, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 2).Value = 278500
$ws.Cells.Item(216, 3).Value = 681153

$ws.Cells.Item(217, 1).Value = 45975
$ws.Cells.Item(217, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 2).Value = 272500
$ws.Cells.Item(217, 3).Value = 579023

# ---- Sheet 2 ----
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(210, 4).Value = 248412
$ws.Cells.Item(210, 5).Value = 393789270
$ws.Cells.Item(210, 6).Value = 27573732000
$ws.Cells.Item(210, 7).Value = 43710608970000
$ws.Cells.Item(210, 8).Value = 0.05999999865889549

$ws.Cells.Item(211, 4).Value = 256318
$ws.Cells.Item(211, 5).Value = 393789270
$ws.Cells.Item(211, 6).Value = 28707616000
$ws.Cells.Item(211, 7).Value = 44104398240000
$ws.Cells.Item(211, 8).Value = 0.07000000029802322

$ws.Cells.Item(212, 1).Value = 45968
$ws.Cells.Item(212, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 2).Value = 110100
$ws.Cells.Item(212, 3).Value = 536664
$ws.Cells.Item(212, 4).Value = 243797
$ws.Cells.Item(212, 5).Value = 393789270
$ws.Cells.Item(212, 6).Value = 26842049700
$ws.Cells.Item(212, 7).Value = 43356198627000
$ws.Cells.Item(212, 8).Value = 0.05999999865889549

$ws.Cells.Item(213, 1).Value = 45971
$ws.Cells.Item(213, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213, 2).Value = 113700
$ws.Cells.Item(213, 3).Value = 699593
$ws.Cells.Item(213, 4).Value = 243320
$ws.Cells.Item(213, 5).Value = 393789270
$ws.Cells.Item(213, 6).Value = 27665484000
$ws.Cells.Item(213, 7).Value = 44773839999000
$ws.Cells.Item(213, 8).Value = 0.05999999865889549

$ws.Cells.Item(214, 1).Value = 45972
$ws.Cells.Item(214, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214, 2).Value = 116000
$ws.Cells.Item(214, 3).Value = 959520
$ws.Cells.Item(214, 4).Value = 251058
$ws.Cells.Item(214, 5).Value = 393789270
$ws.Cells.Item(214, 6).Value = 29122728000
$ws.Cells.Item(214, 7).Value = 45679555320000
$ws.Cells.Item(214, 8).Value = 0.05999999865889549

$ws.Cells.Item(215, 1).Value = 45973
$ws.Cells.Item(215, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 2).Value = 118600
$ws.Cells.Item(215, 3).Value = 1238570
$ws.Cells.Item(215, 4).Value = 250683
$ws.Cells.Item(215, 5).Value = 393789270
$ws.Cells.Item(215, 6).Value = 29731003800
$ws.Cells.Item(215, 7).Value = 46703407422000
$ws.Cells.Item(215, 8).Value = 0.05999999865889549

$ws.Cells.Item(216, 1).Value = 45974
$ws.Cells.Item(216, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 2).Value = 118000
$ws.Cells.Item(216, 3).Value = 783292

$ws.Cells.Item(217, 1).Value = 45975
$ws.Cells.Item(217, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 2).Value = 117000
$ws.Cells.Item(217, 3).Value = 667473

# ---- Sheet 3 ----
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(210, 4).Value = 79710
$ws.Cells.Item(210, 5).Value = 91795094
$ws.Cells.Item(210, 6).Value = 22956480000
$ws.Cells.Item(210, 7).Value = 26436987072000
$ws.Cells.Item(210, 8).Value = 0.09000000357627869

$ws.Cells.Item(211, 4).Value = 50141
$ws.Cells.Item(211, 5).Value = 91795094
$ws.Cells.Item(211, 6).Value = 14766524500
$ws.Cells.Item(211, 7).Value = 27033655183000
$ws.Cells.Item(211, 8).Value = 0.05000000074505806

$ws.Cells.Item(212, 1).Value = 45968
$ws.Cells.Item(212, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 2).Value = 290000
$ws.Cells.Item(212, 3).Value = 127187
$ws.Cells.Item(212, 4).Value = 49127
$ws.Cells.Item(212, 5).Value = 91795094
$ws.Cells.Item(212, 6).Value = 14246830000
$ws.Cells.Item(212, 7).Value = 26620577260000
$ws.Cells.Item(212, 8).Value = 0.05000000074505806

$ws.Cells.Item(213, 1).Value = 45971
$ws.Cells.Item(213, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213, 2).Value = 300500
$ws.Cells.Item(213, 3).Value = 181548
$ws.Cells.Item(213, 4).Value = 48855
$ws.Cells.Item(213, 5).Value = 91795094
$ws.Cells.Item(213, 6).Value = 14680927500
$ws.Cells.Item(213, 7).Value = 27584425747000
$ws.Cells.Item(213, 8).Value = 0.05000000074505806

$ws.Cells.Item(214, 1).Value = 45972
$ws.Cells.Item(214, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214, 2).Value = 298500
$ws.Cells.Item(214, 3).Value = 163489
$ws.Cells.Item(214, 4).Value = 48855
$ws.Cells.Item(214, 5).Value = 91795094
$ws.Cells.Item(214, 6).Value = 14583217500
$ws.Cells.Item(214, 7).Value = 27400835559000
$ws.Cells.Item(214, 8).Value = 0.05000000074505806

$ws.Cells.Item(215, 1).Value = 45973
$ws.Cells.Item(215, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 2).Value = 300000
$ws.Cells.Item(215, 3).Value = 175934
$ws.Cells.Item(215, 4).Value = 47884
$ws.Cells.Item(215, 5).Value = 91795094
$ws.Cells.Item(215, 6).Value = 14365200000
$ws.Cells.Item(215, 7).Value = 27538528200000
$ws.Cells.Item(215, 8).Value = 0.05000000074505806

$ws.Cells.Item(216, 1).Value = 45974
$ws.Cells.Item(216, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 2).Value = 300000
$ws.Cells.Item(216, 3).Value = 167662

$ws.Cells.Item(217, 1).Value = 45975
$ws.Cells.Item(217, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 2).Value = 298000
$ws.Cells.Item(217, 3).Value = 134033

# ---- Sheet 4 ----
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(210, 4).Value = 806497
$ws.Cells.Item(210, 5).Value = 880000000
$ws.Cells.Item(210, 6).Value = 22138342650
$ws.Cells.Item(210, 7).Value = 24156000000000
$ws.Cells.Item(210, 8).Value = 0.09000000357627869

$ws.Cells.Item(211, 3).Value = 8090477
$ws.Cells.Item(211, 4).Value = 808496
$ws.Cells.Item(211, 5).Value = 880000000
$ws.Cells.Item(211, 6).Value = 21142170400
$ws.Cells.Item(211, 7).Value = 23012000000000
$ws.Cells.Item(211, 8).Value = 0.09000000357627869

$ws.Cells.Item(212, 1).Value = 45968
$ws.Cells.Item(212, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 2).Value = 26300
$ws.Cells.Item(212, 3).Value = 7392826
$ws.Cells.Item(212, 4).Value = 805680
$ws.Cells.Item(212, 5).Value = 880000000
$ws.Cells.Item(212, 6).Value = 21189384000
$ws.Cells.Item(212, 7).Value = 23144000000000
$ws.Cells.Item(212, 8).Value = 0.09000000357627869

$ws.Cells.Item(213, 1).Value = 45971
$ws.Cells.Item(213, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213, 2).Value = 27050
$ws.Cells.Item(213, 3).Value = 4105719
$ws.Cells.Item(213, 4).Value = 801820
$ws.Cells.Item(213, 5).Value = 880000000
$ws.Cells.Item(213, 6).Value = 21689231000
$ws.Cells.Item(213, 7).Value = 23804000000000
$ws.Cells.Item(213, 8).Value = 0.09000000357627869

$ws.Cells.Item(214, 1).Value = 45972
$ws.Cells.Item(214, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214, 2).Value = 26400
$ws.Cells.Item(214, 3).Value = 4352613
$ws.Cells.Item(214, 4).Value = 800820
$ws.Cells.Item(214, 5).Value = 880000000
$ws.Cells.Item(214, 6).Value = 21141648000
$ws.Cells.Item(214, 7).Value = 23232000000000
$ws.Cells.Item(214, 8).Value = 0.09000000357627869

$ws.Cells.Item(215, 1).Value = 45973
$ws.Cells.Item(215, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 2).Value = 26250
$ws.Cells.Item(215, 3).Value = 3838511
$ws.Cells.Item(215, 4).Value = 788820
$ws.Cells.Item(215, 5).Value = 880000000
$ws.Cells.Item(215, 6).Value = 20706525000
$ws.Cells.Item(215, 7).Value = 23100000000000
$ws.Cells.Item(215, 8).Value = 0.09000000357627869

$ws.Cells.Item(216, 1).Value = 45974
$ws.Cells.Item(216, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 2).Value = 26500
$ws.Cells.Item(216, 3).Value = 4577984

$ws.Cells.Item(217, 1).Value = 45975
$ws.Cells.Item(217, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 2).Value = 26300
$ws.Cells.Item(217, 3).Value = 7153065

# ---- Sheet 5 ----
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(210, 4).Value = 322550
$ws.Cells.Item(210, 5).Value = 88773116
$ws.Cells.Item(210, 6).Value = 170306400000
$ws.Cells.Item(210, 7).Value = 46872205248000
$ws.Cells.Item(210, 8).Value = 0.3600000143051147

$ws.Cells.Item(211, 3).Value = 161433
$ws.Cells.Item(211, 4).Value = 331901
$ws.Cells.Item(211, 5).Value = 88773116
$ws.Cells.Item(211, 6).Value = 174248025000
$ws.Cells.Item(211, 7).Value = 46605885900000
$ws.Cells.Item(211, 8).Value = 0.3700000047683716

$ws.Cells.Item(212, 1).Value = 45968
$ws.Cells.Item(212, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 2).Value = 521000
$ws.Cells.Item(212, 3).Value = 267709
$ws.Cells.Item(212, 4).Value = 364951
$ws.Cells.Item(212, 5).Value = 88773116
$ws.Cells.Item(212, 6).Value = 190139471000
$ws.Cells.Item(212, 7).Value = 46250793436000
$ws.Cells.Item(212, 8).Value = 0.4099999964237213

$ws.Cells.Item(213, 1).Value = 45971
$ws.Cells.Item(213, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213, 2).Value = 538000
$ws.Cells.Item(213, 3).Value = 182710
$ws.Cells.Item(213, 4).Value = 394493
$ws.Cells.Item(213, 5).Value = 88773116
$ws.Cells.Item(213, 6).Value = 212237234000
$ws.Cells.Item(213, 7).Value = 47759936408000
$ws.Cells.Item(213, 8).Value = 0.4399999976158142

$ws.Cells.Item(214, 1).Value = 45972
$ws.Cells.Item(214, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214, 2).Value = 534000
$ws.Cells.Item(214, 3).Value = 262860
$ws.Cells.Item(214, 4).Value = 466314
$ws.Cells.Item(214, 5).Value = 88773116
$ws.Cells.Item(214, 6).Value = 249011676000
$ws.Cells.Item(214, 7).Value = 47404843944000
$ws.Cells.Item(214, 8).Value = 0.5299999713897705

$ws.Cells.Item(215, 1).Value = 45973
$ws.Cells.Item(215, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 2).Value = 536000
$ws.Cells.Item(215, 3).Value = 231897
$ws.Cells.Item(215, 4).Value = 485961
$ws.Cells.Item(215, 5).Value = 88773116
$ws.Cells.Item(215, 6).Value = 260475096000
$ws.Cells.Item(215, 7).Value = 47582390176000
$ws.Cells.Item(215, 8).Value = 0.550000011920929

$ws.Cells.Item(216, 1).Value = 45974
$ws.Cells.Item(216, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 2).Value = 568000
$ws.Cells.Item(216, 3).Value = 345051

$ws.Cells.Item(217, 1).Value = 45975
$ws.Cells.Item(217, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 2).Value = 586000
$ws.Cells.Item(217, 3).Value = 1117932

# ---- Sheet 6 ----
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(210, 4).Value = 560175
$ws.Cells.Item(210, 5).Value = 306413394
$ws.Cells.Item(210, 6).Value = 70806120000
$ws.Cells.Item(210, 7).Value = 38730653001600
$ws.Cells.Item(210, 8).Value = 0.1800000071525574

$ws.Cells.Item(211, 3).Value = 2066794
$ws.Cells.Item(211, 4).Value = 534199
$ws.Cells.Item(211, 5).Value = 306413394
$ws.Cells.Item(211, 6).Value = 65706477000
$ws.Cells.Item(211, 7).Value = 37688847462000
$ws.Cells.Item(211, 8).Value = 0.1700000017881393

$ws.Cells.Item(212, 1).Value = 45968
$ws.Cells.Item(212, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(212, 2).Value = 126800
$ws.Cells.Item(212, 3).Value = 4418401
$ws.Cells.Item(212, 4).Value = 497865
$ws.Cells.Item(212, 5).Value = 306413394
$ws.Cells.Item(212, 6).Value = 63129282000
$ws.Cells.Item(212, 7).Value = 38853218359200
$ws.Cells.Item(212, 8).Value = 0.1599999964237213

$ws.Cells.Item(213, 1).Value = 45971
$ws.Cells.Item(213, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(213, 2).Value = 129600
$ws.Cells.Item(213, 3).Value = 1860111
$ws.Cells.Item(213, 4).Value = 498973
$ws.Cells.Item(213, 5).Value = 306413394
$ws.Cells.Item(213, 6).Value = 64666900800
$ws.Cells.Item(213, 7).Value = 39711175862400
$ws.Cells.Item(213, 8).Value = 0.1599999964237213

$ws.Cells.Item(214, 1).Value = 45972
$ws.Cells.Item(214, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(214, 2).Value = 126600
$ws.Cells.Item(214, 3).Value = 2065914
$ws.Cells.Item(214, 4).Value = 501442
$ws.Cells.Item(214, 5).Value = 306413394
$ws.Cells.Item(214, 6).Value = 63482557200
$ws.Cells.Item(214, 7).Value = 38791935680400
$ws.Cells.Item(214, 8).Value = 0.1599999964237213

$ws.Cells.Item(215, 1).Value = 45973
$ws.Cells.Item(215, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(215, 2).Value = 126500
$ws.Cells.Item(215, 3).Value = 1031163
$ws.Cells.Item(215, 4).Value = 491785
$ws.Cells.Item(215, 5).Value = 306413394
$ws.Cells.Item(215, 6).Value = 62210802500
$ws.Cells.Item(215, 7).Value = 38761294341000
$ws.Cells.Item(215, 8).Value = 0.1599999964237213

$ws.Cells.Item(216, 1).Value = 45974
$ws.Cells.Item(216, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(216, 2).Value = 130500
$ws.Cells.Item(216, 3).Value = 1627651

$ws.Cells.Item(217, 1).Value = 45975
$ws.Cells.Item(217, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(217, 2).Value = 129100
$ws.Cells.Item(217, 3).Value = 2817706

Write-Output "done"
